$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Update the "score" column (B) on Sheet1 for every student row (2..121)
# from 11.11 to 2020.11
for ($r = 2; $r -le 121; $r++) {
    $ws1.Cells.Item($r, 2).Value = 2020.11
}

# Sheet3 was the active/searched sheet before; leave it parked at B2 with
# no page scrolled (topLeftCell reset) and give it a print page setup.
[void]$ws3.Activate()
[void]$ws3.Range("B2").Select()
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Finish on Sheet1 at B116, which becomes the active/selected tab.
[void]$ws1.Activate()
[void]$ws1.Range("B116").Select()
